$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right below the header row (row 2) and shift existing data down
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the latest dividend entry
$ws.Range("A2").Value = "16/09/2025"
$ws.Range("B2").Value = "16/09/2025"

# Gross Dividend values in this sheet are stored as text (e.g. "0.027"),
# not numbers, so force text formatting before assigning, then restore
# the cell to the sheet's normal style so no stray formatting is visible.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0.027"
$ws.Range("C2").Style = "Normal"
